$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (A1:R1): drop the bold font / thin border / center+top alignment style
# (styles.xml: fonts 2->1, borders 2->1, cellXfs 2->1) so header cells fall back to the default style ---
$ws.Range("A1:R1").ClearFormats()

# --- Refreshed Copernicus metadata: stacOrCswTbox (L) / tempExtentEnd (Q) / numLayers (O) / tempExtentBegin (I) ---
$ws.Range('L2').Value = '[''2020-11-01T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L3').Value = '[''2022-11-01T03:00:00.000Z'', ''2026-02-03T00:00:00.000Z'']'
$ws.Range('L4').Value = '[''2021-10-01T00:00:00.000Z'', ''2026-02-02T00:00:00.000Z'']'
$ws.Range('L10').Value = '[''2021-07-05T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L11').Value = '[''2022-08-01T00:00:00.000Z'', ''2026-02-01T12:00:00.000Z'']'
$ws.Range('L12').Value = '[''2019-01-01T00:00:00.000Z'', ''2026-02-02T00:00:00.000Z'']'
$ws.Range('L13').Value = '[''2019-08-01T00:00:00.000Z'', ''2026-01-31T23:00:00.000Z'']'
$ws.Range('L19').Value = '[''2021-11-01T00:00:00.000Z'', ''2026-02-02T12:00:00.000Z'']'
$ws.Range('L20').Value = '[''2021-10-01T01:00:00.000Z'', ''2026-02-02T12:00:00.000Z'']'
$ws.Range('L21').Value = '[''2022-11-01T00:00:00.000Z'', ''2026-02-01T00:00:00.000Z'']'
$ws.Range('L25').Value = '[''2021-11-01T00:00:00.000Z'', ''2026-02-02T00:00:00.000Z'']'
$ws.Range('L26').Value = '[''2022-12-18T00:00:00.000Z'', ''2026-02-03T11:00:00.000Z'']'
$ws.Range('L27').Value = '[''2023-10-01T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L31').Value = '[''2022-11-22T00:00:00.000Z'', ''2026-02-03T00:00:00.000Z'']'
$ws.Range('L32').Value = '[''2022-11-26T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L33').Value = '[''2022-11-23T00:00:00.000Z'', ''2026-01-31T00:00:00.000Z'']'
$ws.Range('L37').Value = '[''2023-11-01T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L38').Value = '[''2021-11-30T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L39').Value = '[''2023-10-01T00:00:00.000Z'', ''2026-02-02T00:00:00.000Z'']'
$ws.Range('L41').Value = '[''1985-01-01T00:00:00.000Z'', ''2025-12-31T23:00:00.000Z'']'
$ws.Range('L42').Value = '[''1999-01-01T00:00:00.000Z'', ''2025-12-31T00:00:00.000Z'']'
$ws.Range('L43').Value = '[''2022-11-22T00:00:00.000Z'', ''2026-02-03T00:00:00.000Z'']'
$ws.Range('L44').Value = '[''2023-09-29T00:00:00.000Z'', ''2026-01-27T00:00:00.000Z'']'
$ws.Range('L45').Value = '[''2022-11-26T00:00:00.000Z'', ''2026-02-02T23:00:00.000Z'']'
$ws.Range('L46').Value = '[''2022-11-23T00:00:00.000Z'', ''2026-01-31T00:00:00.000Z'']'
$ws.Range('L50').Value = '[''2022-05-03T00:00:00.000Z'', ''2026-01-24T17:44:27.000Z'']'
$ws.Range('L51').Value = '[''2022-10-04T00:00:00.000Z'', ''2026-01-25T00:00:00.000Z'']'
$ws.Range('L55').Value = '[''2022-03-14T01:03:00.000Z'', ''2026-01-24T12:30:42.127Z'']'
$ws.Range('L56').Value = '[''2022-10-04T00:00:00.000Z'', ''2026-01-25T00:00:00.000Z'']'
$ws.Range('L64').Value = '[''2023-04-25T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L65').Value = '[''2023-04-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L70').Value = '[''2023-04-25T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L71').Value = '[''2025-10-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L72').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L74').Value = '[''2023-04-18T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L76').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L78').Value = '[''2023-04-29T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L79').Value = '[''2022-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L80').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L82').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L84').Value = '[''2023-04-29T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L85').Value = '[''2022-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L86').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L90').Value = '[''1997-09-04T00:00:00.000Z'', ''2026-01-17T00:00:00.000Z'']'
$ws.Range('L91').Value = '[''1997-09-01T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L92').Value = '[''1997-09-04T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L94').Value = '[''1997-09-04T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L95').Value = '[''1997-09-01T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L96').Value = '[''1997-09-04T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L98').Value = '[''1997-09-16T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L99').Value = '[''1997-09-01T00:00:00.000Z'', ''2026-01-12T00:00:00.000Z'']'
$ws.Range('L100').Value = '[''1997-09-16T00:00:00.000Z'', ''2026-01-16T00:00:00.000Z'']'
$ws.Range('L101').Value = '[''1997-09-01T00:00:00.000Z'', ''2026-01-12T00:00:00.000Z'']'
$ws.Range('L102').Value = '[''2020-12-20T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L103').Value = '[''2024-01-17T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L104').Value = '[''2020-12-20T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L105').Value = '[''2018-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L106').Value = '[''2019-03-11T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L107').Value = '[''2023-01-30T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L109').Value = '[''2008-01-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L110').Value = '[''2008-01-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L111').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-23T23:00:00.000Z'']'
$ws.Range('L112').Value = '[''2008-01-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L113').Value = '[''2008-01-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L114').Value = '[''2019-01-01T00:00:00.000Z'', ''2026-01-23T23:00:00.000Z'']'
$ws.Range('L122').Value = '[''1982-01-01T00:00:00.000Z'', ''2025-12-25T00:00:00.000Z'']'
$ws.Range('L123').Value = '[''1982-01-01T00:00:00.000Z'', ''2025-12-25T00:00:00.000Z'']'
$ws.Range('L124').Value = '[''1982-01-01T00:00:00.000Z'', ''2025-12-25T00:00:00.000Z'']'
$ws.Range('L125').Value = '[''1982-01-01T00:00:00.000Z'', ''2025-12-25T00:00:00.000Z'']'
$ws.Range('L126').Value = '[''2022-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L127').Value = '[''2023-10-18T00:00:00.000Z'', ''2026-01-19T00:00:00.000Z'']'
$ws.Range('L128').Value = '[''2022-01-02T00:00:00.000Z'', ''2026-01-24T12:00:00.000Z'']'
$ws.Range('L129').Value = '[''2023-02-02T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L130').Value = '[''2025-11-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L131').Value = '[''2018-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L132').Value = '[''2020-09-01T12:10:00.000Z'', ''2026-01-24T10:10:00.000Z'']'
$ws.Range('L133').Value = '[''2023-11-25T06:40:00.000Z'', ''2026-01-24T21:28:00.000Z'']'
$ws.Range('L135').Value = '[''2018-01-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('Q137').Value = '''2025-09-30'
$ws.Range('L143').Value = '[''2024-10-03T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L146').Value = '[''2016-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L147').Value = '[''2023-11-20T00:00:00.000Z'', ''2026-01-23T23:00:00.000Z'']'
$ws.Range('L148').Value = '[''2024-04-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L149').Value = '[''2024-04-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L150').Value = '[''2024-04-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L151').Value = '[''2024-04-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L152').Value = '[''2024-04-01T00:00:00.000Z'', ''2026-01-24T00:00:00.000Z'']'
$ws.Range('L161').Value = '[''2021-01-01T00:00:00.000Z'', ''2026-01-24T21:41:45.000Z'']'
$ws.Range('L162').Value = '[''2021-01-01T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L164').Value = '[''2023-06-08T00:00:00.000Z'', ''2026-01-31T21:00:00.000Z'']'
$ws.Range('L170').Value = '[''1986-06-02T09:00:00.000Z'', ''2026-01-24T07:24:22.000Z'']'
$ws.Range('O171').Value = 838
$ws.Range('L179').Value = '[''1841-03-21T00:00:00.000Z'', ''2026-01-24T22:10:00.000Z'']'
$ws.Range('L180').Value = '[''1841-03-21T00:00:00.000Z'', ''2026-01-21T14:38:00.000Z'']'
$ws.Range('L182').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-24T22:15:00.000Z'']'
$ws.Range('L183').Value = '[''2020-01-01T00:00:00.000Z'', ''2026-01-24T22:15:00.000Z'']'
$ws.Range('L185').Value = '[''2023-11-01T00:00:00.000Z'', ''2026-01-23T19:50:00.000Z'']'
$ws.Range('L186').Value = '[''1993-01-01T00:00:00.000Z'', ''2026-01-18T00:00:00.000Z'']'
$ws.Range('O186').Value = 30
$ws.Range('L187').Value = '[''1993-01-06T00:00:00.000Z'', ''2026-01-14T00:00:00.000Z'']'
$ws.Range('L189').Value = '[''1993-01-01T00:00:00.000Z'', ''2026-01-23T23:00:00.000Z'']'
$ws.Range('L191').Value = '[''2010-01-12T00:00:00.000Z'', ''2026-01-23T00:00:00.000Z'']'
$ws.Range('L193').Value = '[''1993-01-01T00:00:00.000Z'', ''2026-01-18T00:00:00.000Z'']'
$ws.Range('I223').Value = '''1993-01-01'
$ws.Range('L223').Value = '[''1993-01-01T00:00:00.000Z'', ''2023-12-31T00:00:00.000Z'']'
$ws.Range('Q223').Value = '''2023-12-31'
